# Commit: "added test for searchRelatedRecords"
# Adds a new "testBrokerSearchRelatedRecords" test-data block (rows 26-29)
# to the DashboardPageData sheet, mirroring the format of the existing
# "testPoliciesDashboardUI" block right above it (rows 20-23).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DashboardPageData")

# ---------------------------------------------------------------------
# 1. Copy formatting from the nearest existing block so the new rows
#    pick up the same (already-registered) cell styles instead of
#    minting brand-new ones.
# ---------------------------------------------------------------------
$ws.Range("A20").Copy()
$ws.Range("A26").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B20").Copy()
$ws.Range("B26").PasteSpecial(-4122)

$ws.Range("A21").Copy()
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("B21").Copy()
$ws.Range("B27").PasteSpecial(-4122)
$ws.Range("C21").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D21").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E21").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F21").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("G21").Copy()
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("H21").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("F21").Copy()
$ws.Range("I27").PasteSpecial(-4122)
$ws.Range("F21").Copy()
$ws.Range("J27").PasteSpecial(-4122)

$ws.Range("A22").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("B22").Copy()
$ws.Range("B28").PasteSpecial(-4122)
$ws.Range("C22").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D22").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D22").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F22").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("G22").Copy()
$ws.Range("G28").PasteSpecial(-4122)
$ws.Range("H22").Copy()
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range("F22").Copy()
$ws.Range("I28").PasteSpecial(-4122)
$ws.Range("F22").Copy()
$ws.Range("J28").PasteSpecial(-4122)

$ws.Range("A23").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("B23").Copy()
$ws.Range("B29").PasteSpecial(-4122)
$ws.Range("C23").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("D23").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D23").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("F23").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("G23").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("H23").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("F23").Copy()
$ws.Range("I29").PasteSpecial(-4122)
$ws.Range("F23").Copy()
$ws.Range("J29").PasteSpecial(-4122)

$ws.Application.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Values. Shared strings get interned in the order they're assigned,
#    so cells are filled in the same (non-linear) sequence as the
#    original edit to keep the shared-string table in the same order.
# ---------------------------------------------------------------------

# Row 26: section title
$ws.Cells.Item(26, 1).Value = "testBrokerSearchRelatedRecords"

# Numeric / reused-string cells (no new shared strings interned here).
$ws.Cells.Item(27, 1).Value = "runMode"
$ws.Cells.Item(27, 2).Value = "brokerId"
$ws.Cells.Item(27, 3).Value = "agentId"
$ws.Cells.Item(27, 4).Value = "agencyOfficeId"
$ws.Cells.Item(28, 1).Value = "N"
$ws.Cells.Item(28, 2).Value = 20217
$ws.Cells.Item(28, 3).Value = 173
$ws.Cells.Item(28, 4).Value = 237
$ws.Cells.Item(29, 1).Value = "Y"
$ws.Cells.Item(29, 2).Value = 20217
$ws.Cells.Item(29, 3).Value = 237
$ws.Cells.Item(29, 4).Value = 8006

# referenceNumber number formats: set N-row (General) before Y-row (Text)
# so the two new cell styles register in the same order as the target
# (style index 9 = General/N-row, style index 10 = Text/Y-row).
$ws.Cells.Item(28, 5).NumberFormat = "General"
$ws.Cells.Item(29, 5).NumberFormat = "@"

# referenceNumber header
$ws.Cells.Item(27, 5).Value = "referenceNumber"
# quoteName value (Y row) typed ahead of its own header
$ws.Cells.Item(29, 6).Value = "Test Partner API"
# quoteName / policyName headers
$ws.Cells.Item(27, 6).Value = "quoteName"
$ws.Cells.Item(27, 7).Value = "policyName"
# policyName value (Y row)
$ws.Cells.Item(29, 7).Value = "Kelley Buick GMC"
# referenceNumber value (Y row) - kept as text
$ws.Cells.Item(29, 5).Value = "12825076"
# noSuchARecord value (Y row)
$ws.Cells.Item(29, 8).Value = "afd98afd"
# noSuchARecord / expForNoSuchARecord headers
$ws.Cells.Item(27, 8).Value = "noSuchARecord"
$ws.Cells.Item(27, 9).Value = "expForNoSuchARecord"
# expForNoSuchARecord value (Y row)
$ws.Cells.Item(29, 9).Value = "Your search has no results"
# policyNumber header
$ws.Cells.Item(27, 10).Value = "policyNumber"
# policyNumber value (Y row)
$ws.Cells.Item(29, 10).Value = "H20NPP70603-00"

# referenceNumber value (N row) - literal number, not a shared string.
$ws.Cells.Item(28, 5).Value = 12825076

# ---------------------------------------------------------------------
# 3. View state: scrolled down with the last populated cell selected.
# ---------------------------------------------------------------------
$ws.Application.Goto($ws.Range("B7"))
$ws.Range("I29").Select()
